# =========================================================================
# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" worksheet (fund holdings detail) positioned
#    right before the "总计" (grand-total) worksheet.
# 2) Prepend a new "2022-Q1" summary row to the "总计" worksheet.
# =========================================================================

$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# 1) New "2022-Q1" sheet
# -------------------------------------------------------------------------
# Copy "2021-Q4" (it already has the right header/column styling) and drop
# the copy immediately before "总计"; re-fetch "总计" by name afterwards
# since sheet handles here are positional, not identity-based.
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$q4Sheet.Copy($totalSheet, $null)

$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# "2021-Q4" only had 9 data rows (rows 2-10); "2022-Q1" needs 12 data rows
# (rows 2-13) so extend the already-styled row 10 down by 3 more rows
# (format first, then values, so the new rows keep the bordered/bold style
# on column A instead of silently reverting to the default style).
$newSheet.Range("A10:H10").Copy()
$newSheet.Range("A11:H13").PasteSpecial(-4122)
$newSheet.Range("A10:H10").Copy()
$newSheet.Range("A11:H13").PasteSpecial(-4163)

$data2022Q1 = New-Object 'object[,]' 12,8
$data2022Q1[0,0] = 0
$data2022Q1[0,1] = '''002345'
$data2022Q1[0,2] = '华夏高端制造灵活配置混合'
$data2022Q1[0,3] = '''23.47'
$data2022Q1[0,4] = '''93.20'
$data2022Q1[0,5] = '''3.67'
$data2022Q1[0,6] = '''0.8613'
$data2022Q1[0,7] = 8
$data2022Q1[1,0] = 1
$data2022Q1[1,1] = '''009861'
$data2022Q1[1,2] = '鹏华新兴成长混合A'
$data2022Q1[1,3] = '''23.92'
$data2022Q1[1,4] = '''90.71'
$data2022Q1[1,5] = '''3.53'
$data2022Q1[1,6] = '''0.8444'
$data2022Q1[1,7] = 8
$data2022Q1[2,0] = 2
$data2022Q1[2,1] = '''005028'
$data2022Q1[2,2] = '鹏华研究精选灵活配置混合'
$data2022Q1[2,3] = '''15.55'
$data2022Q1[2,4] = '''87.74'
$data2022Q1[2,5] = '''3.31'
$data2022Q1[2,6] = '''0.5147'
$data2022Q1[2,7] = 9
$data2022Q1[3,0] = 3
$data2022Q1[3,1] = '''004640'
$data2022Q1[3,2] = '华夏节能环保股票'
$data2022Q1[3,3] = '''5.45'
$data2022Q1[3,4] = '''90.73'
$data2022Q1[3,5] = '''5.13'
$data2022Q1[3,6] = '''0.2796'
$data2022Q1[3,7] = 5
$data2022Q1[4,0] = 4
$data2022Q1[4,1] = '''004558'
$data2022Q1[4,2] = '汇安丰裕灵活配置混合A'
$data2022Q1[4,3] = '''0.99'
$data2022Q1[4,4] = '''83.41'
$data2022Q1[4,5] = '''3.99'
$data2022Q1[4,6] = '''0.0395'
$data2022Q1[4,7] = 3
$data2022Q1[5,0] = 5
$data2022Q1[5,1] = '''011603'
$data2022Q1[5,2] = '兴业高端制造混合A'
$data2022Q1[5,3] = '''1.19'
$data2022Q1[5,4] = '''76.40'
$data2022Q1[5,5] = '''2.73'
$data2022Q1[5,6] = '''0.0325'
$data2022Q1[5,7] = 7
$data2022Q1[6,0] = 6
$data2022Q1[6,1] = '''009862'
$data2022Q1[6,2] = '鹏华新兴成长混合C'
$data2022Q1[6,3] = '''0.73'
$data2022Q1[6,4] = '''90.71'
$data2022Q1[6,5] = '''3.53'
$data2022Q1[6,6] = '''0.0258'
$data2022Q1[6,7] = 8
$data2022Q1[7,0] = 7
$data2022Q1[7,1] = '''011604'
$data2022Q1[7,2] = '兴业高端制造混合C'
$data2022Q1[7,3] = '''0.54'
$data2022Q1[7,4] = '''76.40'
$data2022Q1[7,5] = '''2.73'
$data2022Q1[7,6] = '''0.0147'
$data2022Q1[7,7] = 7
$data2022Q1[8,0] = 8
$data2022Q1[8,1] = '''005437'
$data2022Q1[8,2] = '易方达易百智能量化策略灵活配置混合A'
$data2022Q1[8,3] = '''0.96'
$data2022Q1[8,4] = '''94.52'
$data2022Q1[8,5] = '''1.08'
$data2022Q1[8,6] = '''0.0104'
$data2022Q1[8,7] = 6
$data2022Q1[9,0] = 9
$data2022Q1[9,1] = '''005438'
$data2022Q1[9,2] = '易方达易百智能量化策略灵活配置混合C'
$data2022Q1[9,3] = '''0.17'
$data2022Q1[9,4] = '''94.52'
$data2022Q1[9,5] = '''1.08'
$data2022Q1[9,6] = '''0.0018'
$data2022Q1[9,7] = 6
$data2022Q1[10,0] = 10
$data2022Q1[10,1] = '''004559'
$data2022Q1[10,2] = '汇安丰裕灵活配置混合C'
$data2022Q1[10,3] = '''0.01'
$data2022Q1[10,4] = '''83.41'
$data2022Q1[10,5] = '''3.99'
$data2022Q1[10,6] = '''0.0004'
$data2022Q1[10,7] = 3
$data2022Q1[11,0] = 11
$data2022Q1[11,1] = '''003739'
$data2022Q1[11,2] = '新华鑫弘灵活配置混合'
$data2022Q1[11,3] = '''0.01'
$data2022Q1[11,4] = '''42.81'
$data2022Q1[11,5] = '''0.53'
$data2022Q1[11,6] = '''0.0001'
$data2022Q1[11,7] = 7

$newSheet.Range("A2:H13").Value = $data2022Q1

# -------------------------------------------------------------------------
# 2) Prepend the "2022-Q1" row to the "总计" sheet
# -------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Shift existing rows down one (row3->row4, row2->row3), bottom-up so we
# never overwrite data before it has been copied. Two-pass paste (formats,
# then values) keeps column A's bordered/bold style intact.
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A4:D4").PasteSpecial(-4122)
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A4:D4").PasteSpecial(-4163)

$totalSheet.Range("A2:D2").Copy()
$totalSheet.Range("A3:D3").PasteSpecial(-4122)
$totalSheet.Range("A2:D2").Copy()
$totalSheet.Range("A3:D3").PasteSpecial(-4163)

# Re-index column A (0, 1, 2, ...)
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2

# Write the new 2022-Q1 summary row
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 12
$totalSheet.Range("D2").Value = 2.63
